$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 3
$ws.Range("C5").Value = 202

$ws.Range("J5").Select()
